$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet keeps Coin/Link as plain text and Price/Volume(1h) as numeric-looking
# text (e.g. "305.93", "-3.87%"). A leading apostrophe is used only for the
# numeric-looking D/E values so Excel keeps them as text instead of converting
# them to a number/percentage; the plain B/C text needs no such prefix.
$ws.Range("D2").Value = "'305.93"
$ws.Range("E2").Value = "'-3.87%"

$ws.Range("D3").Value = "'37.18"
$ws.Range("E3").Value = "'-6.15%"

$ws.Range("D4").Value = "'5.083"
$ws.Range("E4").Value = "'-1.02%"

$ws.Range("D5").Value = "'0.07709"
$ws.Range("E5").Value = "'-6.06%"

$ws.Range("D6").Value = "'4.351"
$ws.Range("E6").Value = "'1.11%"

$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.907"
$ws.Range("E7").Value = "'-6.88%"

$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'8.200"
$ws.Range("E8").Value = "'-1.88%"

$ws.Range("E9").Value = "'-3.45%"

$ws.Range("D10").Value = "'0.9164"
$ws.Range("E10").Value = "'-2.21%"

$ws.Range("D11").Value = "'0.1175"
$ws.Range("E11").Value = "'-13.15%"

$ws.Range("D12").Value = "'0.1866"
$ws.Range("E12").Value = "'-6.32%"

$ws.Range("D13").Value = "'0.08669"
$ws.Range("E13").Value = "'-4.38%"

$ws.Range("D14").Value = "'0.03404"
$ws.Range("E14").Value = "'-3.00%"

$ws.Range("D15").Value = "'0.09687"
$ws.Range("E15").Value = "'-1.12%"

$ws.Range("D16").Value = "'0.001365"
$ws.Range("E16").Value = "'-2.82%"

$ws.Range("D17").Value = "'0.005921"
$ws.Range("E17").Value = "'-5.50%"

$ws.Range("D18").Value = "'3.613"
$ws.Range("E18").Value = "'-1.89%"

$ws.Range("D19").Value = "'0.3410"
$ws.Range("E19").Value = "'-2.03%"

$ws.Range("D20").Value = "'0.1275"
$ws.Range("E20").Value = "'-3.55%"

$ws.Range("D21").Value = "'5.015"
$ws.Range("E21").Value = "'0.23%"

$ws.Range("D22").Value = "'0.2595"
$ws.Range("E22").Value = "'5.95%"

$ws.Range("E23").Value = "'5,173.77%"

$ws.Range("D24").Value = "'0.04326"

$ws.Range("D25").Value = "'0.001213"
$ws.Range("E25").Value = "'-1.21%"

$ws.Range("D26").Value = "'0.004536"
$ws.Range("E26").Value = "'-5.60%"

$ws.Range("D27").Value = "'0.0001352"
$ws.Range("E27").Value = "'3.84%"

$ws.Range("D39").Value = "'0.02200"
$ws.Range("E39").Value = "'-2.77%"

$ws.Range("D40").Value = "'0.04905"
$ws.Range("E40").Value = "'-5.47%"

$ws.Range("D41").Value = "'0.007534"
$ws.Range("E41").Value = "'-2.74%"

$ws.Range("D42").Value = "'0.009870"
$ws.Range("E42").Value = "'0.14%"

$ws.Range("D43").Value = "'0.1330"
$ws.Range("E43").Value = "'-5.46%"

$ws.Range("D44").Value = "'0.001996"
$ws.Range("E44").Value = "'-2.39%"

$ws.Range("D45").Value = "'0.008491"
$ws.Range("E45").Value = "'-12.18%"

$ws.Range("D46").Value = "'0.00006545"
$ws.Range("E46").Value = "'-0.80%"

$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.19%"

$ws.Range("D48").Value = "'0.002998"
$ws.Range("E48").Value = "'1.83%"

$ws.Range("D49").Value = "'0.001302"
$ws.Range("E49").Value = "'-22.93%"

$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.19%"

$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.19%"
